# correctif scrapping functions 3SN1 3NI1
#
# Appends freshly-scraped rows to the "2M30", "3NI1" and "3SN1" price
# sheets, then refreshes the "RPA" summary sheet: updates the latest
# NI/SN quotes and drops the stale ZN line (ZN is no longer scraped here).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 2M30 : three more scraped rows (18/10/2023) — price column failed to
# scrape for these, so only date / currency / unit came back.
# ---------------------------------------------------------------------
$ws2M30 = $wb.Worksheets.Item("2M30")

$ws2M30.Range("A117").Value = "18/10/2023"
$ws2M30.Range("C117").Value = "€"
$ws2M30.Range("D117").Value = "100KG"

$ws2M30.Range("A118").Value = "18/10/2023"
$ws2M30.Range("C118").Value = "€"
$ws2M30.Range("D118").Value = "100KG"

$ws2M30.Range("A119").Value = "18/10/2023"
$ws2M30.Range("C119").Value = "€"
$ws2M30.Range("D119").Value = "100KG"

# ---------------------------------------------------------------------
# 3NI1 : row 94 completes 18/10/2023, rows 95-97 are scraper misfires
# (currency/unit only, no date or price), rows 98-99 bring 19/10/2023.
# ---------------------------------------------------------------------
$ws3NI1 = $wb.Worksheets.Item("3NI1")

$ws3NI1.Range("A94").Value = "18/10/2023"
$ws3NI1.Range("B94").Value = "18.470,00"
$ws3NI1.Range("C94").Value = "$"
$ws3NI1.Range("D94").Value = "TO"

$ws3NI1.Range("C95").Value = "$"
$ws3NI1.Range("D95").Value = "TO"

$ws3NI1.Range("C96").Value = "$"
$ws3NI1.Range("D96").Value = "TO"

$ws3NI1.Range("C97").Value = "$"
$ws3NI1.Range("D97").Value = "TO"

$ws3NI1.Range("A98").Value = "19/10/2023"
$ws3NI1.Range("B98").Value = "18.285,00"
$ws3NI1.Range("C98").Value = "$"
$ws3NI1.Range("D98").Value = "TO"

$ws3NI1.Range("A99").Value = "19/10/2023"
$ws3NI1.Range("B99").Value = "18.285,00"
$ws3NI1.Range("C99").Value = "$"
$ws3NI1.Range("D99").Value = "TO"

# ---------------------------------------------------------------------
# 3SN1 : rows 98-102 finish out 18/10/2023, row 103 brings 19/10/2023.
# ---------------------------------------------------------------------
$ws3SN1 = $wb.Worksheets.Item("3SN1")

$ws3SN1.Range("A98").Value = "18/10/2023"
$ws3SN1.Range("B98").Value = "25.450,00"
$ws3SN1.Range("C98").Value = "$"
$ws3SN1.Range("D98").Value = "TO"

$ws3SN1.Range("A99").Value = "18/10/2023"
$ws3SN1.Range("B99").Value = "25.450,00"
$ws3SN1.Range("C99").Value = "$"
$ws3SN1.Range("D99").Value = "TO"

$ws3SN1.Range("A100").Value = "18/10/2023"
$ws3SN1.Range("B100").Value = "25.450,00"
$ws3SN1.Range("C100").Value = "$"
$ws3SN1.Range("D100").Value = "TO"

$ws3SN1.Range("A101").Value = "18/10/2023"
$ws3SN1.Range("B101").Value = "25.450,00"
$ws3SN1.Range("C101").Value = "$"
$ws3SN1.Range("D101").Value = "TO"

$ws3SN1.Range("A102").Value = "18/10/2023"
$ws3SN1.Range("B102").Value = "25.450,00"
$ws3SN1.Range("C102").Value = "$"
$ws3SN1.Range("D102").Value = "TO"

$ws3SN1.Range("A103").Value = "19/10/2023"
$ws3SN1.Range("B103").Value = "25.075,00"
$ws3SN1.Range("C103").Value = "$"
$ws3SN1.Range("D103").Value = "TO"

# ---------------------------------------------------------------------
# RPA : refresh the latest NI / SN quotes, and remove the ZN row (row 4)
# entirely — it is no longer part of the scraped summary.
# ---------------------------------------------------------------------
$wsRPA = $wb.Worksheets.Item("RPA")

$wsRPA.Range("C2").Value = "18.285,00"
$wsRPA.Range("C3").Value = "25.075,00"
$wsRPA.Rows.Item(4).Delete()
